# Insert a new "Skill Description" column after column A (SkillCode),
# duplicating the SkillCode value into the new column for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (SkillCode)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# Insert a new column before column B, shifting SFIA Level/Keycode/Description right.
$ws.Columns.Item(2).Insert()

# New header for the inserted column
$ws.Cells.Item(1, 2).Value = "Skill Description"

# Fill the new column with the SkillCode text for each data row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value()
}
